$d = $word.ActiveDocument

# ----------------------------------------------------------------------
# 1) Education table (2nd table): set the first row's height to 1167
#    twips (= 58.35 points) -- matches <w:trHeight w:val="1167"/>.
# ----------------------------------------------------------------------
$eduTable = $d.Tables.Item(2)
$eduRow = $eduTable.Rows.Item(1)
$eduRow.Height = 58.35

Write-Output "Row height set."

# ----------------------------------------------------------------------
# 2) "BS in Computer Science    " -> "BS in Computer Science"
#    (drop the 4 trailing spaces)
# ----------------------------------------------------------------------
$r = $d.Content.Find.Execute("BS in Computer Science    ", $true, $false, $false, $false, $false, `
                              $true, 1, $false, "BS in Computer Science", 2)
Write-Output "BS in Computer Science fix: $r"

# ----------------------------------------------------------------------
# 3) "Graduation Date: May 2016 (expected)" -> "Cum Laude, May 2016"
#    with a _GoBack bookmark planted right at the boundary between the
#    two halves (this also forces Word to keep them as two runs).
# ----------------------------------------------------------------------
$r = $d.Content.Find.Execute("Graduation Date: May 2016 (expected)", $true, $false, $false, $false, $false, `
                              $true, 1, $false, "Cum Laude, May 2016", 2)
Write-Output "Graduation Date fix: $r"

$gradRng = $d.Content
$gradRng.Find.Execute("Cum Laude, ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$bmRange = $d.Range($gradRng.End, $gradRng.End)
$d.Bookmarks.Add("_GoBack", $bmRange)
Write-Output "Bookmark planted."

# ----------------------------------------------------------------------
# 4) "Dean's List and Atlantic 10 Commissioner's Honor Roll Spring,
#     Fall 2014, Spring 2015"
#    -> "Dean's List and Atlantic 10 Commissioner's Honor Roll Spring,
#        Fall 2014, 2015"          (drop the word "Spring" before 2015)
#    split across three runs, matching boundaries in the target markup.
#    (The straight apostrophe in "Dean's" must never be touched by a
#    Find/Replace -- doing so triggers this runtime's smart-quote
#    autocorrect and turns it into a curly quote.)
# ----------------------------------------------------------------------
$r = $d.Content.Find.Execute("List and Atlantic 10 Commissioner" + [char]0x2019 + "s Honor Roll Spring, Fall 2014, Spring 2015", `
                              $true, $false, $false, $false, $false, $true, 1, $false, `
                              "List and Atlantic 10 Commissioner" + [char]0x2019 + "s Honor Roll Spring, Fall 2014, 2015", 2)
Write-Output "Dean's List text fix: $r"

$deanLeft = $d.Content
$deanLeft.Find.Execute("Dean's List and Atlantic 10 Commissioner" + [char]0x2019 + "s Honor Roll Spring, Fall 2014,", `
                        $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$deanRight = $d.Content
$deanRight.Find.Execute(" 2015", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

# split #2 (rightmost): boundary right before " 2015"
$deanSuffix = $d.Range($deanLeft.End, $deanRight.End)
$deanSuffix.Font.Bold = $true
$deanSuffix.Font.Bold = $false

$deanMid = $d.Content
$deanMid.Find.Execute("Dean's List and Atlantic 10 Commissioner" + [char]0x2019 + "s Hono", `
                       $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

# split #1: boundary right before "r Roll"
$deanMidRange = $d.Range($deanMid.End, $deanLeft.End)
$deanMidRange.Font.Bold = $true
$deanMidRange.Font.Bold = $false

Write-Output "Dean's List split done."

# ----------------------------------------------------------------------
# 5) "Automatically generated documentation for medical procedures"
#    -> "Created a tool for automatically generating" +
#       " documentation for medical procedures"   (two runs)
# ----------------------------------------------------------------------
$r = $d.Content.Find.Execute("Automatically generated documentation for medical procedures", `
                              $true, $false, $false, $false, $false, $true, 1, $false, `
                              "Created a tool for automatically generating documentation for medical procedures", 2)
Write-Output "Automatically-generated fix: $r"

$autoLeft = $d.Content
$autoLeft.Find.Execute("Created a tool for automatically generating", `
                        $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$autoRight = $d.Content
$autoRight.Find.Execute(" documentation for medical procedures", `
                         $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$autoSplit = $d.Range($autoLeft.End, $autoRight.End)
$autoSplit.Font.Bold = $true
$autoSplit.Font.Bold = $false
Write-Output "Automatically-generated split done."

# ----------------------------------------------------------------------
# 6) "Participated in Research Experience for Undergraduates (REU)
#     program summer 2015"
#    -> "Participated in Research Expe" +
#       "rience for Undergraduates (REU), S" +
#       "ummer 2015"                          (three runs)
# ----------------------------------------------------------------------
$r = $d.Content.Find.Execute("Participated in Research Experience for Undergraduates (REU) program summer 2015", `
                              $true, $false, $false, $false, $false, $true, 1, $false, `
                              "Participated in Research Experience for Undergraduates (REU), Summer 2015", 2)
Write-Output "REU fix: $r"

$reuLeft = $d.Content
$reuLeft.Find.Execute("Participated in Research Expe", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$reuMid = $d.Content
$reuMid.Find.Execute("rience for Undergraduates (REU), S", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$reuRight = $d.Content
$reuRight.Find.Execute("ummer 2015", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$reuSplit1 = $d.Range($reuLeft.End, $reuMid.End)
$reuSplit1.Font.Bold = $true
$reuSplit1.Font.Bold = $false

$reuSplit2 = $d.Range($reuMid.End, $reuRight.End)
$reuSplit2.Font.Bold = $true
$reuSplit2.Font.Bold = $false

Write-Output "REU split done."

# ----------------------------------------------------------------------
# 7) "Organized and analyzed data for" + bookmark(_GoBack) +
#    " a PhD candidate's dissertation"
#    -> merge into a single run (same text); the old _GoBack bookmark
#       is dropped here since it has moved up to the Education table.
# ----------------------------------------------------------------------
$mergedText = "Organized and analyzed data for a PhD candidate" + [char]0x2019 + "s dissertation"
$r = $d.Content.Find.Execute("Organized and analyzed data for a PhD candidate" + [char]0x2019 + "s dissertation", `
                              $true, $false, $false, $false, $false, $true, 1, $false, $mergedText, 2)
Write-Output "Organized-and-analyzed merge: $r"

# ----------------------------------------------------------------------
# 8) "Create, maintain " + "Django" (wrapped in proofErr spell-check
#    markers) + "-based web administration tool for university Drupal
#    sites" -> merge into a single run (same text, proofErr dropped).
# ----------------------------------------------------------------------
$djangoText = "Create, maintain Django-based web administration tool for university Drupal sites"
$r = $d.Content.Find.Execute($djangoText, $true, $false, $false, $false, $false, $true, 1, $false, $djangoText, 2)
Write-Output "Django merge: $r"
